# Update column G ("K" - strikeouts) values per regenerated save_data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 0
    3  = 2
    4  = 0
    5  = 0
    6  = 1
    7  = 0
    8  = 0
    9  = 1
    10 = 0
    11 = 0
    12 = 2
    13 = 1
    14 = 1
    15 = 1
    16 = 0
    17 = 2
    18 = 2
    19 = 0
    20 = 1
    21 = 2
    22 = 1
    23 = 1
    24 = 0
    25 = 1
    26 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
